# Apply update: add data for 2022-04-07
# - Rename sheet "Through 2022-03-29" -> "Through 2022-03-30"
# - Update header cell I1 text "2022 (through 03-29)" -> "2022 (through 03-30)"
# - Update March value I4: 125 -> 130
# - Update Total value I14: 425 -> 430

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet tab
$ws.Name = "Through 2022-03-30"

# Update the header label in I1 (shared string text)
$ws.Range("I1").Value = "2022 (through 03-30)"

# Update March (row 4) value for column I
$ws.Range("I4").Value = 130

# Update Total (row 14) value for column I
$ws.Range("I14").Value = 430
